# Applies the re-ordering / update of the species-observation rows (Artfynd sheet):
#   - existing rows 2-6 are re-sorted into a new order (data values unchanged)
#   - a brand-new observation (Spillkraka, row 7) is appended
# The simplest reliable way to express a full row re-shuffle is to clear the data
# rows and rewrite them in the desired final order/content in one bulk range write.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove any pre-existing data rows below the header so stale cells from the old
# row order do not linger once the new row order/content is written.
$ws.Rows("2:100").Clear()

# A few columns hold values that *look* numeric/date-like but are plain text in the
# source data (e.g. Startdatum "2017-10-14", Antal "1"). Pre-format those columns as
# Text so Excel does not "helpfully" convert them into date serials / numbers on entry.
$ws.Range("I6:I7").NumberFormat = "@"
$ws.Range("Y2:Y7").NumberFormat = "@"
$ws.Range("AA2:AA7").NumberFormat = "@"

# cols: A, B, C, D, E, F, G, H, I, J, K, L, M, N, O, P, Q, R, S, T, U, V, W, X, Y, Z, AA, AB, AC, AD, AE, AF, AG, AH, AI, AJ, AK, AL, AM, AN, AO, AP, AQ, AR, AS, AT, AU, AV, AW, AX, AY
$data = New-Object 'object[,]' 6,51
# Row 2
$data[0,0] = 68146859  # A2
$data[0,1] = 90674  # B2
$data[0,2] = 'Ovaliderad'  # C2
$data[0,3] = 'LC'  # D2
$data[0,4] = 5964  # E2
$data[0,5] = 'Fjällig taggsvamp s.str.'  # F2
$data[0,6] = 'Sarcodon imbricatus s.str.'  # G2
$data[0,7] = '(L.:Fr.) P.Karst.'  # H2
$data[0,8] = $null  # I2
$data[0,9] = $null  # J2
$data[0,10] = $null  # K2
$data[0,11] = $null  # L2
$data[0,12] = $null  # M2
$data[0,13] = $null  # N2
$data[0,14] = $null  # O2
$data[0,15] = 'Nygård, 1,1 km S-ut, Upl'  # P2
$data[0,16] = 653808.9988882283  # Q2
$data[0,17] = 6643999.173554195  # R2
$data[0,18] = 10  # S2
$data[0,19] = 'Uppsala'  # T2
$data[0,20] = 'Uppsala'  # U2
$data[0,21] = 'Uppland'  # V2
$data[0,22] = 'Vaksala'  # W2
$data[0,23] = $null  # X2
$data[0,24] = '2017-10-14'  # Y2
$data[0,25] = '00:00'  # Z2
$data[0,26] = '2017-10-14'  # AA2
$data[0,27] = '00:00'  # AB2
$data[0,28] = $null  # AC2
$data[0,29] = $false  # AD2
$data[0,30] = $false  # AE2
$data[0,31] = $null  # AF2
$data[0,32] = $false  # AG2
$data[0,33] = $null  # AH2
$data[0,34] = $null  # AI2
$data[0,35] = $null  # AJ2
$data[0,36] = $null  # AK2
$data[0,37] = $null  # AL2
$data[0,38] = $null  # AM2
$data[0,39] = $null  # AN2
$data[0,40] = $null  # AO2
$data[0,41] = $null  # AP2
$data[0,42] = $null  # AQ2
$data[0,43] = $null  # AR2
$data[0,44] = $null  # AS2
$data[0,45] = $null  # AT2
$data[0,46] = $null  # AU2
$data[0,47] = $null  # AV2
$data[0,48] = 'Henry Åkerström'  # AW2
$data[0,49] = 'Henry Åkerström, Thorleif Joelson'  # AX2
$data[0,50] = 'Naturskyddsföreningen Uppsala, skogsgruppexkursion'  # AY2

# Row 3
$data[1,0] = 68145535  # A3
$data[1,1] = 98520  # B3
$data[1,2] = 'Ovaliderad'  # C3
$data[1,3] = 'LC'  # D3
$data[1,4] = 222498  # E3
$data[1,5] = 'Blåsippa'  # F3
$data[1,6] = 'Hepatica nobilis'  # G3
$data[1,7] = 'Schreb.'  # H3
$data[1,8] = $null  # I3
$data[1,9] = $null  # J3
$data[1,10] = $null  # K3
$data[1,11] = $null  # L3
$data[1,12] = $null  # M3
$data[1,13] = $null  # N3
$data[1,14] = $null  # O3
$data[1,15] = 'Nygård, 1,0 km S-ut, Upl'  # P3
$data[1,16] = 653895.8525253957  # Q3
$data[1,17] = 6644088.162086425  # R3
$data[1,18] = 10  # S3
$data[1,19] = 'Uppsala'  # T3
$data[1,20] = 'Uppsala'  # U3
$data[1,21] = 'Uppland'  # V3
$data[1,22] = 'Vaksala'  # W3
$data[1,23] = $null  # X3
$data[1,24] = '2017-10-14'  # Y3
$data[1,25] = '00:00'  # Z3
$data[1,26] = '2017-10-14'  # AA3
$data[1,27] = '00:00'  # AB3
$data[1,28] = $null  # AC3
$data[1,29] = $false  # AD3
$data[1,30] = $false  # AE3
$data[1,31] = $null  # AF3
$data[1,32] = $false  # AG3
$data[1,33] = $null  # AH3
$data[1,34] = $null  # AI3
$data[1,35] = $null  # AJ3
$data[1,36] = $null  # AK3
$data[1,37] = $null  # AL3
$data[1,38] = $null  # AM3
$data[1,39] = $null  # AN3
$data[1,40] = $null  # AO3
$data[1,41] = $null  # AP3
$data[1,42] = $null  # AQ3
$data[1,43] = $null  # AR3
$data[1,44] = $null  # AS3
$data[1,45] = $null  # AT3
$data[1,46] = $null  # AU3
$data[1,47] = $null  # AV3
$data[1,48] = 'Henry Åkerström'  # AW3
$data[1,49] = 'Henry Åkerström, Thorleif Joelson'  # AX3
$data[1,50] = 'Naturskyddsföreningen Uppsala, skogsgruppexkursion'  # AY3

# Row 4
$data[2,0] = 68146871  # A4
$data[2,1] = 4711  # B4
$data[2,2] = 'Ovaliderad'  # C4
$data[2,3] = 'LC'  # D4
$data[2,4] = 100299  # E4
$data[2,5] = 'Thomsons trägnagare'  # F4
$data[2,6] = 'Cacotemnus thomsoni'  # G4
$data[2,7] = '(Kraatz, 1881)'  # H4
$data[2,8] = $null  # I4
$data[2,9] = $null  # J4
$data[2,10] = $null  # K4
$data[2,11] = $null  # L4
$data[2,12] = 'äldre gnagspår'  # M4
$data[2,13] = $null  # N4
$data[2,14] = $null  # O4
$data[2,15] = 'Nygård, 1,0 km S-ut, Upl'  # P4
$data[2,16] = 653842.0029573618  # Q4
$data[2,17] = 6644028.167799042  # R4
$data[2,18] = 10  # S4
$data[2,19] = 'Uppsala'  # T4
$data[2,20] = 'Uppsala'  # U4
$data[2,21] = 'Uppland'  # V4
$data[2,22] = 'Vaksala'  # W4
$data[2,23] = $null  # X4
$data[2,24] = '2017-10-14'  # Y4
$data[2,25] = '00:00'  # Z4
$data[2,26] = '2017-10-14'  # AA4
$data[2,27] = '00:00'  # AB4
$data[2,28] = $null  # AC4
$data[2,29] = $false  # AD4
$data[2,30] = $false  # AE4
$data[2,31] = $null  # AF4
$data[2,32] = $false  # AG4
$data[2,33] = $null  # AH4
$data[2,34] = $null  # AI4
$data[2,35] = 'gran'  # AJ4
$data[2,36] = 'Picea abies'  # AK4
$data[2,37] = $null  # AL4
$data[2,38] = 'Stående död trädstam/högstubbe'  # AM4
$data[2,39] = $null  # AN4
$data[2,40] = 'Standing dead tree/snags # Picea abies'  # AO4
$data[2,41] = $null  # AP4
$data[2,42] = $null  # AQ4
$data[2,43] = $null  # AR4
$data[2,44] = $null  # AS4
$data[2,45] = $null  # AT4
$data[2,46] = $null  # AU4
$data[2,47] = $null  # AV4
$data[2,48] = 'Henry Åkerström'  # AW4
$data[2,49] = 'Henry Åkerström, Thorleif Joelson'  # AX4
$data[2,50] = 'Naturskyddsföreningen Uppsala, skogsgruppexkursion'  # AY4

# Row 5
$data[3,0] = 68146872  # A5
$data[3,1] = 5113  # B5
$data[3,2] = 'Ovaliderad'  # C5
$data[3,3] = 'LC'  # D5
$data[3,4] = 100526  # E5
$data[3,5] = 'Bronshjon'  # F5
$data[3,6] = 'Callidium coriaceum'  # G5
$data[3,7] = 'Paykull, 1800'  # H5
$data[3,8] = $null  # I5
$data[3,9] = $null  # J5
$data[3,10] = $null  # K5
$data[3,11] = $null  # L5
$data[3,12] = 'äldre gnagspår'  # M5
$data[3,13] = $null  # N5
$data[3,14] = $null  # O5
$data[3,15] = 'Nygård, 1,0 km S-ut, Upl'  # P5
$data[3,16] = 653842.0029573618  # Q5
$data[3,17] = 6644028.167799042  # R5
$data[3,18] = 10  # S5
$data[3,19] = 'Uppsala'  # T5
$data[3,20] = 'Uppsala'  # U5
$data[3,21] = 'Uppland'  # V5
$data[3,22] = 'Vaksala'  # W5
$data[3,23] = $null  # X5
$data[3,24] = '2017-10-14'  # Y5
$data[3,25] = '00:00'  # Z5
$data[3,26] = '2017-10-14'  # AA5
$data[3,27] = '00:00'  # AB5
$data[3,28] = $null  # AC5
$data[3,29] = $false  # AD5
$data[3,30] = $false  # AE5
$data[3,31] = $null  # AF5
$data[3,32] = $false  # AG5
$data[3,33] = $null  # AH5
$data[3,34] = $null  # AI5
$data[3,35] = 'gran'  # AJ5
$data[3,36] = 'Picea abies'  # AK5
$data[3,37] = $null  # AL5
$data[3,38] = 'Stående död trädstam/högstubbe'  # AM5
$data[3,39] = $null  # AN5
$data[3,40] = 'Standing dead tree/snags # Picea abies'  # AO5
$data[3,41] = $null  # AP5
$data[3,42] = $null  # AQ5
$data[3,43] = $null  # AR5
$data[3,44] = $null  # AS5
$data[3,45] = $null  # AT5
$data[3,46] = $null  # AU5
$data[3,47] = $null  # AV5
$data[3,48] = 'Henry Åkerström'  # AW5
$data[3,49] = 'Henry Åkerström, Thorleif Joelson'  # AX5
$data[3,50] = 'Naturskyddsföreningen Uppsala, skogsgruppexkursion'  # AY5

# Row 6
$data[4,0] = 109439564  # A6
$data[4,1] = 56717  # B6
$data[4,2] = 'Ovaliderad'  # C6
$data[4,3] = 'NT'  # D6
$data[4,4] = 103008  # E6
$data[4,5] = 'Ärtsångare'  # F6
$data[4,6] = 'Curruca curruca'  # G6
$data[4,7] = '(Linnaeus, 1758)'  # H6
$data[4,8] = '1'  # I6
$data[4,9] = $null  # J6
$data[4,10] = $null  # K6
$data[4,11] = $null  # L6
$data[4,12] = 'spel/sång'  # M6
$data[4,13] = $null  # N6
$data[4,14] = $null  # O6
$data[4,15] = 'Skogsängen, Upl'  # P6
$data[4,16] = 653910.4788543681  # Q6
$data[4,17] = 6644087.263707791  # R6
$data[4,18] = 25  # S6
$data[4,19] = 'Uppsala'  # T6
$data[4,20] = 'Uppsala'  # U6
$data[4,21] = 'Uppland'  # V6
$data[4,22] = 'Vaksala'  # W6
$data[4,23] = $null  # X6
$data[4,24] = '2023-05-23'  # Y6
$data[4,25] = '11:03'  # Z6
$data[4,26] = '2023-05-23'  # AA6
$data[4,27] = '11:03'  # AB6
$data[4,28] = $null  # AC6
$data[4,29] = $false  # AD6
$data[4,30] = $false  # AE6
$data[4,31] = $null  # AF6
$data[4,32] = $false  # AG6
$data[4,33] = $null  # AH6
$data[4,34] = $null  # AI6
$data[4,35] = $null  # AJ6
$data[4,36] = $null  # AK6
$data[4,37] = $null  # AL6
$data[4,38] = $null  # AM6
$data[4,39] = $null  # AN6
$data[4,40] = $null  # AO6
$data[4,41] = $null  # AP6
$data[4,42] = $null  # AQ6
$data[4,43] = $null  # AR6
$data[4,44] = $null  # AS6
$data[4,45] = $null  # AT6
$data[4,46] = $null  # AU6
$data[4,47] = $null  # AV6
$data[4,48] = 'Jimmy Peterson'  # AW6
$data[4,49] = 'Jimmy Peterson'  # AX6
$data[4,50] = $null  # AY6

# Row 7
$data[5,0] = 112536523  # A7
$data[5,1] = 56446  # B7
$data[5,2] = 'Ovaliderad'  # C7
$data[5,3] = 'NT'  # D7
$data[5,4] = 100049  # E7
$data[5,5] = 'Spillkråka'  # F7
$data[5,6] = 'Dryocopus martius'  # G7
$data[5,7] = '(Linnaeus, 1758)'  # H7
$data[5,8] = '1'  # I7
$data[5,9] = $null  # J7
$data[5,10] = $null  # K7
$data[5,11] = $null  # L7
$data[5,12] = 'lockläte, övriga läten'  # M7
$data[5,13] = $null  # N7
$data[5,14] = $null  # O7
$data[5,15] = 'Skogsängen, Upl'  # P7
$data[5,16] = 653910  # Q7
$data[5,17] = 6644087  # R7
$data[5,18] = 25  # S7
$data[5,19] = 'Uppsala'  # T7
$data[5,20] = 'Uppsala'  # U7
$data[5,21] = 'Uppland'  # V7
$data[5,22] = 'Vaksala'  # W7
$data[5,23] = $null  # X7
$data[5,24] = '2023-10-05'  # Y7
$data[5,25] = '15:50'  # Z7
$data[5,26] = '2023-10-05'  # AA7
$data[5,27] = '15:50'  # AB7
$data[5,28] = $null  # AC7
$data[5,29] = $false  # AD7
$data[5,30] = $false  # AE7
$data[5,31] = $null  # AF7
$data[5,32] = $false  # AG7
$data[5,33] = $null  # AH7
$data[5,34] = $null  # AI7
$data[5,35] = $null  # AJ7
$data[5,36] = $null  # AK7
$data[5,37] = $null  # AL7
$data[5,38] = $null  # AM7
$data[5,39] = $null  # AN7
$data[5,40] = $null  # AO7
$data[5,41] = $null  # AP7
$data[5,42] = $null  # AQ7
$data[5,43] = $null  # AR7
$data[5,44] = $null  # AS7
$data[5,45] = $null  # AT7
$data[5,46] = $null  # AU7
$data[5,47] = $null  # AV7
$data[5,48] = 'Jimmy Peterson'  # AW7
$data[5,49] = 'Jimmy Peterson'  # AX7
$data[5,50] = $null  # AY7

$ws.Range("A2:AY7").Value = $data

$ws.Range("A1").Select()